$wb = $excel.ActiveWorkbook

# Data for the new row (row 27) in each of the four worksheets.
$rowsData = @{
    "DE_LFT_#1" = @{
        A = 45813.43701388889
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x78"
        E = "0x14"
        F = 380
        G = 759863127514710900000000.0
        H = 376
        I = 14
    }
    "DE_LFT_#2" = @{
        A = 45813.43701388889
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x78"
        E = "0xe"
        F = 380
        G = 568432987514711000000000.0
        H = 376
        I = 14
    }
    "DE_PLT_#1" = @{
        A = 45813.43701388889
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x82"
        E = "0x7"
        F = 130
        G = 568631262647114000000000.0
        H = 130
        I = 7
    }
    "DE_PLT_#2" = @{
        A = 45813.43701388889
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x81"
        E = "0x3"
        F = 130
        G = 985046333984776000000000.0
        H = 129
        I = 3
    }
}

foreach ($sheetName in $rowsData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $rowsData[$sheetName]

    $ws.Cells.Item(27, 1).Value = $data.A
    $ws.Cells.Item(27, 1).NumberFormat = $ws.Cells.Item(26, 1).NumberFormat
    $ws.Cells.Item(27, 2).Value = $data.B
    $ws.Cells.Item(27, 3).Value = $data.C
    $ws.Cells.Item(27, 4).Value = $data.D
    $ws.Cells.Item(27, 5).Value = $data.E
    $ws.Cells.Item(27, 6).Value = $data.F
    $ws.Cells.Item(27, 7).Value = $data.G
    $ws.Cells.Item(27, 8).Value = $data.H
    $ws.Cells.Item(27, 9).Value = $data.I
}
